# Apply BOM updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add C20, C21 to capacitor designator list (1uF, 603, qty stays 8)
$ws.Range("A4").Value = "C1, C11, C12, C13, C14, C2, C20, C21, C5, C9"

# Row 11: add D2 to Schottky diode designator list
$ws.Range("A11").Value = "D17, D18, D2"

# Row 20: replace Q4 (BSC098N10NS5 MOSFET, TDSON-8-1) with Q3 (Si4435DDY-T1-E3, SOIC-8)
$ws.Range("A20").Value = "Q3"
$ws.Range("B20").Value = "SOIC127P600X175-8N"
$ws.Range("D20").Value = "Si4435DDY-T1-E3"
$ws.Range("E20").Value = "781-SI4435DDY-T1-E3"
$ws.Range("F20").Value = 2.54
$ws.Range("G20").Value = 2.28
$ws.Range("H20").Value = 1.72

# Row 22: 10k resistor qty 6 -> 8
$ws.Range("C22").Value = 8

# Row 26: R19 (200) qty 1 -> 2
$ws.Range("C26").Value = 2

# Row 39: SW25 -> SW3, SW25 ; SW_RST -> Switch
$ws.Range("A39").Value = "SW3, SW25"
$ws.Range("D39").Value = "Switch"

# Row 30: R22, R23 56k -> R22, R23, R26, R27 5.1k (new part/footprint/pricing)
$ws.Range("A30").Value = "R22, R23, R26, R27"
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = "5.1k"
$ws.Range("E30").Value = "603-RT0805DRE075K1L"
$ws.Range("F30").Value = 0.56
$ws.Range("G30").Value = 0.472
$ws.Range("H30").Value = 0.12

# Style update: B4:B8 footprint column gets numeric-left-aligned style
$ws.Range("B4:B8").NumberFormat = "0"
$ws.Range("B4:B8").HorizontalAlignment = -4131

# Sheet view: clear frozen/scrolled topLeftCell, move selection to C17
$ws.Range("A1").Select()
$ws.Range("C17").Select()
